# LA Clippers roster update:
#  - Remove "Moses Brown (TW)" (was row 10) from the roster.
#  - Re-sequence the "No." index column (A) for the rows that shift up.
#  - Fill in previously-missing jersey numbers (No. column B) for Eric Gordon,
#    Bones Hyland and Mason Plumlee, and re-order Hyland/Plumlee.
#  - Add a new player, "Nate Darling (TW)", as the final roster row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Moses Brown (TW) row entirely; this shifts rows 11-17 up to 10-16
# (and keeps the existing hyperlink relationships untouched, same as Excel's
# native row-delete behaviour).
$ws.Rows.Item(10).Delete()

# Fix up the sequential index in column A for the rows that moved.
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14

# Eric Gordon now has a jersey number on file.
$ws.Range("B14").Value = 10

# Bones Hyland and Mason Plumlee swap places (rows 15/16) and both now have
# jersey numbers on file.
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = "Bones Hyland"
$ws.Range("D15").Value = "PG"
$ws.Range("E15").Value = "6-3"
$ws.Range("F15").Value = 173
$ws.Range("G15").Value = "September 14, 2000"
$ws.Range("H15").Value = "us"
$ws.Range("I15").Value = "1"
$ws.Range("J15").Value = "VCU"
$ws.Range("K15").Value = "https://www.basketball-reference.com/players/h/hylanbo01.html"

$ws.Range("B16").Value = 44
$ws.Range("C16").Value = "Mason Plumlee"
$ws.Range("D16").Value = "C"
$ws.Range("E16").Value = "6-11"
$ws.Range("F16").Value = 254
$ws.Range("G16").Value = "March 5, 1990"
$ws.Range("H16").Value = "us"
$ws.Range("I16").Value = "9"
$ws.Range("J16").Value = "Duke"
$ws.Range("K16").Value = "https://www.basketball-reference.com/players/p/plumlma01.html"

# Append the new player, Nate Darling (TW), as row 17. Copy formatting from
# row 16 first (bold/bordered index style in A, hyperlink style in K), then
# overwrite with his actual data.
$ws.Range("A16").Copy($ws.Range("A17"))
$ws.Range("K16").Copy($ws.Range("K17"))

$ws.Range("A17").Value = 15
$ws.Range("C17").Value = "Nate Darling (TW)"
$ws.Range("D17").Value = "SG"
$ws.Range("E17").Value = "6-6"
$ws.Range("F17").Value = 200
$ws.Range("G17").Value = "August 30, 1998"
$ws.Range("H17").Value = "ca"
$ws.Range("I17").Value = "1"
$ws.Range("J17").Value = "UAB, University of Delaware"
$ws.Range("K17").Value = "https://www.basketball-reference.com/players/d/darlina01.html"

Write-Output "done"
